$wb = $excel.ActiveWorkbook

function Set-LeveProfitCells {
    param($sheet, $row)
    $eVal = $sheet.Cells.Item($row, 5).Value2
    $kVal = $sheet.Cells.Item($row, 11).Value2
    $lVal = $sheet.Cells.Item($row, 12).Value2
    if ($null -eq $kVal) { $kVal = 0 }
    if ($null -eq $lVal) { $lVal = 0 }
    if ($kVal -ne 0) {
        $sheet.Cells.Item($row, 13).Value = $eVal - $kVal
    } else {
        $sheet.Cells.Item($row, 13).Value = ""
    }
    if ($lVal -ne 0) {
        $sheet.Cells.Item($row, 14).Value = (-2 * $eVal) - $lVal
    } else {
        $sheet.Cells.Item($row, 14).Value = ""
    }
}

$ws = $wb.Worksheets.Item("ALC")

$ws.Cells.Item(16, 8).Value = 24000
$ws.Cells.Item(16, 9).Value = 24000
$ws.Cells.Item(16, 11).Value = 24000
Set-LeveProfitCells $ws 16

$ws.Cells.Item(29, 8).Value = 4837
$ws.Cells.Item(29, 9).Value = 186.66667
$ws.Cells.Item(29, 11).Value = 560.00001
Set-LeveProfitCells $ws 29

$ws.Cells.Item(31, 8).Value = 3678.7856
$ws.Cells.Item(31, 9).Value = 2773.0908
$ws.Cells.Item(31, 11).Value = 8319.2724
Set-LeveProfitCells $ws 31

$ws.Cells.Item(53, 8).Value = 521.5238000000001
$ws.Cells.Item(53, 9).Value = 143.27272
$ws.Cells.Item(53, 10).Value = 937.6
$ws.Cells.Item(53, 11).Value = 143.27272
$ws.Cells.Item(53, 12).Value = 937.6
Set-LeveProfitCells $ws 53

$ws.Cells.Item(100, 8).Value = 11658.167
$ws.Cells.Item(100, 9).Value = 1450
$ws.Cells.Item(100, 10).Value = 13699.8
$ws.Cells.Item(100, 11).Value = 1450
$ws.Cells.Item(100, 12).Value = 13699.8
Set-LeveProfitCells $ws 100

$ws.Cells.Item(111, 8).Value = 1466.6666
$ws.Cells.Item(111, 9).Value = 1450
$ws.Cells.Item(111, 10).Value = 1500
$ws.Cells.Item(111, 11).Value = 4350
$ws.Cells.Item(111, 12).Value = 4500
Set-LeveProfitCells $ws 111

$ws.Cells.Item(131, 8).Value = 3840.2942
$ws.Cells.Item(131, 9).Value = 1553.8462
$ws.Cells.Item(131, 10).Value = 11271.25
$ws.Cells.Item(131, 11).Value = 4661.5386
$ws.Cells.Item(131, 12).Value = 33813.75
Set-LeveProfitCells $ws 131

$ws.Cells.Item(135, 8).Value = 5061.1333
$ws.Cells.Item(135, 9).Value = 6128.381
$ws.Cells.Item(135, 11).Value = 55155.429
Set-LeveProfitCells $ws 135

$ws.Cells.Item(138, 8).Value = 6816.628
$ws.Cells.Item(138, 9).Value = 13260.875
$ws.Cells.Item(138, 11).Value = 39782.625
Set-LeveProfitCells $ws 138


$ws = $wb.Worksheets.Item("ARM")

$ws.Cells.Item(61, 8).Value = 2360
$ws.Cells.Item(61, 9).Value = 2115
$ws.Cells.Item(61, 10).Value = 3993.3333
$ws.Cells.Item(61, 11).Value = 2115
$ws.Cells.Item(61, 12).Value = 3993.3333
Set-LeveProfitCells $ws 61

$ws.Cells.Item(74, 8).Value = 1600.4872
$ws.Cells.Item(74, 9).Value = 1271.8823
$ws.Cells.Item(74, 11).Value = 1271.8823
Set-LeveProfitCells $ws 74

$ws.Cells.Item(77, 8).Value = 1600.4872
$ws.Cells.Item(77, 9).Value = 1271.8823
$ws.Cells.Item(77, 11).Value = 6359.4115
Set-LeveProfitCells $ws 77

$ws.Cells.Item(97, 8).Value = 2588.8235
$ws.Cells.Item(97, 9).Value = 1386.619
$ws.Cells.Item(97, 11).Value = 1386.619
Set-LeveProfitCells $ws 97

$ws.Cells.Item(110, 8).Value = 2263.1428
$ws.Cells.Item(110, 9).Value = 2118.4
$ws.Cells.Item(110, 11).Value = 2118.4
Set-LeveProfitCells $ws 110

$ws.Cells.Item(129, 8).Value = 69799
$ws.Cells.Item(129, 10).Value = 69799
$ws.Cells.Item(129, 12).Value = 69799
Set-LeveProfitCells $ws 129

$ws.Cells.Item(132, 8).Value = 3033.3784
$ws.Cells.Item(132, 9).Value = 3006.5278
$ws.Cells.Item(132, 10).Value = 4000
$ws.Cells.Item(132, 11).Value = 9019.5834
$ws.Cells.Item(132, 12).Value = 12000
Set-LeveProfitCells $ws 132

$ws.Cells.Item(136, 8).Value = 2360
$ws.Cells.Item(136, 9).Value = 2115
$ws.Cells.Item(136, 10).Value = 3993.3333
$ws.Cells.Item(136, 11).Value = 6345
$ws.Cells.Item(136, 12).Value = 11979.9999
Set-LeveProfitCells $ws 136

$ws.Cells.Item(139, 8).Value = 69999.875
$ws.Cells.Item(139, 10).Value = 69999.875
$ws.Cells.Item(139, 12).Value = 69999.875
Set-LeveProfitCells $ws 139


$ws = $wb.Worksheets.Item("BSM")

$ws.Cells.Item(99, 8).Value = 2982.1667
$ws.Cells.Item(99, 9).Value = 2676.4
$ws.Cells.Item(99, 11).Value = 2676.4
Set-LeveProfitCells $ws 99

$ws.Cells.Item(107, 8).Value = 1517.7
$ws.Cells.Item(107, 9).Value = 1517.7
$ws.Cells.Item(107, 11).Value = 1517.7
Set-LeveProfitCells $ws 107

$ws.Cells.Item(134, 8).Value = 5298.028
$ws.Cells.Item(134, 9).Value = 4667.1724
$ws.Cells.Item(134, 11).Value = 14001.5172
Set-LeveProfitCells $ws 134

$ws.Cells.Item(140, 8).Value = 119996
$ws.Cells.Item(140, 10).Value = 119996
$ws.Cells.Item(140, 12).Value = 119996
Set-LeveProfitCells $ws 140


$ws = $wb.Worksheets.Item("CUL")

$ws.Cells.Item(60, 8).Value = 2099.6667
$ws.Cells.Item(60, 9).Value = 2099.6667
$ws.Cells.Item(60, 11).Value = 6299.000100000001
Set-LeveProfitCells $ws 60

$ws.Cells.Item(70, 8).Value = 9580.909
$ws.Cells.Item(70, 9).Value = 1990
$ws.Cells.Item(70, 10).Value = 10340
$ws.Cells.Item(70, 11).Value = 5970
$ws.Cells.Item(70, 12).Value = 31020
Set-LeveProfitCells $ws 70

$ws.Cells.Item(73, 8).Value = 9580.909
$ws.Cells.Item(73, 9).Value = 1990
$ws.Cells.Item(73, 10).Value = 10340
$ws.Cells.Item(73, 11).Value = 5970
$ws.Cells.Item(73, 12).Value = 31020
Set-LeveProfitCells $ws 73

$ws.Cells.Item(86, 8).Value = 434.75
$ws.Cells.Item(86, 9).Value = 413.16666
$ws.Cells.Item(86, 11).Value = 1239.49998
Set-LeveProfitCells $ws 86

$ws.Cells.Item(89, 8).Value = 434.75
$ws.Cells.Item(89, 9).Value = 413.16666
$ws.Cells.Item(89, 11).Value = 3718.49994
Set-LeveProfitCells $ws 89

$ws.Cells.Item(122, 8).Value = 3094.8667
$ws.Cells.Item(122, 10).Value = 3814.9
$ws.Cells.Item(122, 12).Value = 34334.1
Set-LeveProfitCells $ws 122

$ws.Cells.Item(132, 8).Value = 6203.727
$ws.Cells.Item(132, 9).Value = 4379.4
$ws.Cells.Item(132, 11).Value = 39414.6
Set-LeveProfitCells $ws 132

$ws.Cells.Item(139, 8).Value = 18411
$ws.Cells.Item(139, 9).Value = 0
$ws.Cells.Item(139, 11).Value = 0
Set-LeveProfitCells $ws 139


$ws = $wb.Worksheets.Item("GSM")

$ws.Cells.Item(70, 8).Value = 7498.448
$ws.Cells.Item(70, 9).Value = 7658.2856
$ws.Cells.Item(70, 11).Value = 7658.2856
Set-LeveProfitCells $ws 70

$ws.Cells.Item(73, 8).Value = 7498.448
$ws.Cells.Item(73, 9).Value = 7658.2856
$ws.Cells.Item(73, 11).Value = 7658.2856
Set-LeveProfitCells $ws 73

$ws.Cells.Item(86, 8).Value = 44999.75
$ws.Cells.Item(86, 9).Value = 40000
$ws.Cells.Item(86, 11).Value = 40000
Set-LeveProfitCells $ws 86

$ws.Cells.Item(89, 8).Value = 44999.75
$ws.Cells.Item(89, 9).Value = 40000
$ws.Cells.Item(89, 11).Value = 120000
Set-LeveProfitCells $ws 89

$ws.Cells.Item(132, 8).Value = 11550.714
$ws.Cells.Item(132, 9).Value = 0
$ws.Cells.Item(132, 11).Value = 0
Set-LeveProfitCells $ws 132


$ws = $wb.Worksheets.Item("LTW")

$ws.Cells.Item(22, 8).Value = 1181.6428
$ws.Cells.Item(22, 9).Value = 900.7
$ws.Cells.Item(22, 11).Value = 900.7
Set-LeveProfitCells $ws 22

$ws.Cells.Item(26, 8).Value = 30010
$ws.Cells.Item(26, 10).Value = 30010
$ws.Cells.Item(26, 12).Value = 30010
Set-LeveProfitCells $ws 26

$ws.Cells.Item(27, 8).Value = 1181.6428
$ws.Cells.Item(27, 9).Value = 900.7
$ws.Cells.Item(27, 11).Value = 900.7
Set-LeveProfitCells $ws 27

$ws.Cells.Item(31, 8).Value = 14409
$ws.Cells.Item(31, 9).Value = 12947.5
$ws.Cells.Item(31, 10).Value = 17332
$ws.Cells.Item(31, 11).Value = 12947.5
$ws.Cells.Item(31, 12).Value = 17332
Set-LeveProfitCells $ws 31

$ws.Cells.Item(93, 8).Value = 2479
$ws.Cells.Item(93, 9).Value = 2392.3333
$ws.Cells.Item(93, 10).Value = 2999
$ws.Cells.Item(93, 11).Value = 2392.3333
$ws.Cells.Item(93, 12).Value = 2999
Set-LeveProfitCells $ws 93

$ws.Cells.Item(100, 8).Value = 2043.75
$ws.Cells.Item(100, 9).Value = 1173
$ws.Cells.Item(100, 10).Value = 2334
$ws.Cells.Item(100, 11).Value = 1173
$ws.Cells.Item(100, 12).Value = 2334
Set-LeveProfitCells $ws 100

$ws.Cells.Item(122, 8).Value = 3250.2
$ws.Cells.Item(122, 10).Value = 2833
$ws.Cells.Item(122, 12).Value = 8499
Set-LeveProfitCells $ws 122

$ws.Cells.Item(135, 8).Value = 114978.86
$ws.Cells.Item(135, 10).Value = 114978.86
$ws.Cells.Item(135, 12).Value = 114978.86
Set-LeveProfitCells $ws 135

$ws.Cells.Item(136, 8).Value = 3384.5818
$ws.Cells.Item(136, 9).Value = 3300.3076
$ws.Cells.Item(136, 11).Value = 9900.9228
Set-LeveProfitCells $ws 136


$ws = $wb.Worksheets.Item("WVR")

$ws.Cells.Item(113, 8).Value = 664.6316
$ws.Cells.Item(113, 10).Value = 998.2727
$ws.Cells.Item(113, 12).Value = 2994.8181
Set-LeveProfitCells $ws 113

$ws.Cells.Item(122, 8).Value = 47307.91
$ws.Cells.Item(122, 9).Value = 57359.832
$ws.Cells.Item(122, 10).Value = 2074.25
$ws.Cells.Item(122, 11).Value = 172079.496
$ws.Cells.Item(122, 12).Value = 6222.75
Set-LeveProfitCells $ws 122

$ws.Cells.Item(132, 8).Value = 7906.3
$ws.Cells.Item(132, 9).Value = 8444.875
$ws.Cells.Item(132, 11).Value = 25334.625
Set-LeveProfitCells $ws 132

$ws.Cells.Item(136, 8).Value = 7050.6
$ws.Cells.Item(136, 9).Value = 6920.436
$ws.Cells.Item(136, 10).Value = 7896.6665
$ws.Cells.Item(136, 11).Value = 20761.308
$ws.Cells.Item(136, 12).Value = 23689.9995
Set-LeveProfitCells $ws 136

